# Follow channel testcases added
#
# Adds a new "Follow Channels" worksheet (after the existing "Search
# Channels" sheet) containing the UserName / Password / ChannelName /
# runMode test-data row used by the new Follow-Channel test cases.

$wb = $excel.ActiveWorkbook

# The sheet whose header/data row formatting we want to clone, and after
# which the new sheet should be inserted.
$srcSheet = $wb.Worksheets.Item("Search Channels")

# Add the new worksheet right after "Search Channels" (i.e. as the last
# sheet in the workbook) and name it.
$newSheet = $wb.Worksheets.Add($null, $srcSheet)
$newSheet.Name = "Follow Channels"

# Clone the header-row (row 1) and data-row (row 2) cell formatting from
# the "Search Channels" sheet so the new sheet matches the look of the
# other test-data sheets (bold yellow header, bordered data row, etc).
$srcSheet.Range("A1:D1").Copy()
$newSheet.Range("A1:D1").PasteSpecial(-4122)

$srcSheet.Range("A2:D2").Copy()
$newSheet.Range("A2:D2").PasteSpecial(-4122)

# Header row
$newSheet.Range("A1").Value = "UserName"
$newSheet.Range("B1").Value = "Password"
$newSheet.Range("C1").Value = "ChannelName"
$newSheet.Range("D1").Value = "runMode"

# Test-data row
$newSheet.Range("A2").Value = "manisha.kisan17@gmail.com"
$newSheet.Range("B2").Value = "manisha123"
$newSheet.Range("C2").Value = "FCI Magazine"
$newSheet.Range("D2").Value = "Y"
